# "Generate Report for Archive"
# - Status of the two sample rows moves from "Ready for handoff" to
#   "In Translation" on every sheet that surfaces it (Overview's zh-cn/de-de
#   status columns, plus the per-locale "Status" column on the zh-cn and
#   de-de detail sheets).
# - The zh-cn/de-de status columns are narrowed to match the new,
#   shorter heading/value.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Update status text everywhere it appears ---------------------------
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the zh-cn / de-de status columns -----------------------------
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
